# Applies the 'Updated cryptos list' data refresh to Sheet1.
# Only cell VALUES change (prices in column D, volume-% in column E,
# plus a data swap of two coin rows - B/C/D/E for rows 29<->30 and 32<->33).
# Four D-column cells (D9, D26, D27, D46) are forced to Text number-format
# before assignment so Excel doesn't silently drop significant trailing/
# leading zeros (e.g. '1.00' -> 1, '0.0000120' -> 0.000012) when the
# string looks like a plain number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.338.83'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '3.560.92'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D6").Value = '144.68'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").Value = '3.560.22'
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("E9").Value = '  +2.17%  '
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").Value = '7.82'
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '4.164.98'
$ws.Range("D15").Value = '30.31'
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").Value = '3.561.37'
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").Value = '66.404.43'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").Value = '11.47'
$ws.Range("E19").Value = '  +4.62%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '14.83'
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("D22").Value = '431.63'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '0.611'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").Value = '79.49'
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").Value = '3.700.78'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.96'
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '9.14'
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '25.45'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").Value = '3.555.28'
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("E35").Value = '  -6.14%  '
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  -1.34%  '
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").Value = '175.99'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("D41").Value = '0.0847'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("D43").Value = '0.889'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").Value = '45.99'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '2.55'
$ws.Range("E47").Value = '  +5.00%  '
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").Value = '25.29'
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("D50").Value = '7.15'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").Value = '23.23'
$ws.Range("E51").Value = '  +2.82%  '
